# Updates the "cryptos" price list to the latest scraped snapshot.
# Column D (Price) holds values that look numeric ("1.00", "60.817.81", etc.)
# but must stay plain text (as in the source data, dotted big numbers like
# "60.817.81" are not valid numbers anyway). Excel auto-converts a numeric-
# looking string assigned via .Value, so for column D we momentarily force
# the cell to Text format, assign the literal string, then restore the
# original ("Normal") style so no formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.817.81'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.18%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.872.50'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.60%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.09%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.48'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.45%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.98'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.12%  '

# Row 7
$ws.Range('E7').Value = '  +0.02%  '

# Row 8
$ws.Range('E8').Value = '  -3.50%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.77'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.52%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.136'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.97%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.429'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.53%  '

# Row 12
$ws.Range('E12').Value = '  -3.94%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '32.20'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.13%  '

# Row 14
$ws.Range('E14').Value = '  -0.09%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.355.60'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.38%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '60.828.15'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.02%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.881.51'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.24%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.49'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.14%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '423.88'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.48%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.05%  '

# Row 21
$ws.Range('E21').Value = '  -3.90%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.90'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.25%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '79.50'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.35%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.35'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.30%  '

# Row 25
$ws.Range('E25').Value = '  +0.01%  '

# Row 26
$ws.Range('B26').Value = 'Fetch.AI'
$ws.Range('C26').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.03'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -7.25%  '

# Row 27
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.29'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.80%  '

# Row 28
$ws.Range('E28').Value = '  -3.25%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.05'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -10.25%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.67'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.87%  '

# Row 31
$ws.Range('E31').Value = '  +0.12%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '25.52'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.17%  '

# Row 33
$ws.Range('E33').Value = '  -4.65%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0₃0843'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.12%  '

# Row 35
$ws.Range('E35').Value = '  -5.01%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.42'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.74%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '48.90'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.03%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.78'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -8.01%  '

# Row 39
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.117'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.77%  '

# Row 40
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.89'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.64%  '

# Row 41
$ws.Range('E41').Value = '  -2.61%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '38.59'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.63%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.263'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.73%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.653.79'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.43%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '132.71'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.51%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0330'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.80%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '341.82'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -10.11%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '22.35'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.11%  '

# Row 50
$ws.Range('E50').Value = '  -4.19%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.92'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.06%  '
